$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New RFM / derived-feature rows appended to the data dictionary (rows 10-20)
$data = @(
    @("Recency", "Integer", "Measures how recently a customer made a purchase. The more recent the transaction, the more likely the customer is to respond to future promotions.", "No"),
    @("Frequency", "Integer", "Measures how often a customer makes a purchase within a given time period. Frequent buyers are generally more engaged and loyal.", "No"),
    @("Monetary", "Float", "Measures how much a customer spends. High-spending customers are more valuable and should be treated accordingly.", "No"),
    @("TotalQuantity", "Integer", "TotalPrice = Quantity * UnitPrice", "No"),
    @("AvgQuantity", "Float", "average quantity of purchases", "No"),
    @("AvgSpend", "Float", "average spending of customers", "No"),
    @("ProductDiversity", "Integer", "displays the diversity of the products that are purchased", "No"),
    @("R_Score", "Integer", "scores 1-5 for Recency", "No"),
    @("F_Score", "Integer", "scores 1-5 for Frequency", "No"),
    @("M_Score", "Integer", "scores 1-5 for Monetary", "No"),
    @("RFM_Score", "Integer", "scores 1-5 for sum of Recency, Frequency and Monetary", "No")
)

$row = 10
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $row++
}
